$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

$wsOpen   = $wb.Worksheets.Item("Open Tickets")
$wsSolved = $wb.Worksheets.Item("Solved tickets in a year")
$wsOffene = $wb.Worksheets.Item("Offene Tickets")
$wsGesch  = $wb.Worksheets.Item("Geschlossene Tickets pro Jahr")

# ---------------------------------------------------------------------------
# 1) New ticket-category rows on the two "open tickets" sheets.
#    Written in this precise interleaved order so that newly-introduced
#    shared strings are appended to the shared-string table in the same
#    sequence the original authors produced.
# ---------------------------------------------------------------------------

# Offene Tickets (German) rows 10-19 (note: row 13 is entered before row 12)
$wsOffene.Cells.Item(10, 1).Value = "Verbindung per ODBC fehlgeschlagen"
$wsOffene.Cells.Item(10, 2).Value = 300
$wsOffene.Cells.Item(11, 1).Value = "Allgemeine technische Frage"
$wsOffene.Cells.Item(11, 2).Value = 120
$wsOffene.Cells.Item(13, 1).Value = "Keine SAP Daten"
$wsOffene.Cells.Item(13, 2).Value = 200
$wsOffene.Cells.Item(12, 1).Value = "Fehler beim Editieren einer Variable"
$wsOffene.Cells.Item(12, 2).Value = 450
$wsOffene.Cells.Item(14, 1).Value = "Scripting Problem"
$wsOffene.Cells.Item(14, 2).Value = 500
$wsOffene.Cells.Item(15, 1).Value = "Probleme mit dem Bearbeiten von Schriftarten"
$wsOffene.Cells.Item(15, 2).Value = 250
$wsOffene.Cells.Item(16, 1).Value = "Authentifizierung mit JSON nicht möglich"
$wsOffene.Cells.Item(16, 2).Value = 300
$wsOffene.Cells.Item(17, 1).Value = "Download Fehler "
$wsOffene.Cells.Item(17, 2).Value = 300
$wsOffene.Cells.Item(18, 1).Value = "Sonderzeichen im Passwort nicht akzeptiert"
$wsOffene.Cells.Item(18, 2).Value = 100
$wsOffene.Cells.Item(19, 1).Value = "Fehlermeldung"
$wsOffene.Cells.Item(19, 2).Value = 200

# Open Tickets (English) row 20 first
$wsOpen.Cells.Item(20, 1).Value = "Databinding Listview not working"
$wsOpen.Cells.Item(20, 2).Value = 250

# Row 21 - identical text used by both language sheets
$wsOpen.Cells.Item(21, 1).Value = "SQL Async"
$wsOpen.Cells.Item(21, 2).Value = 100

# Open Tickets (English) rows 10-19
$wsOpen.Cells.Item(10, 1).Value = "Connection via ODBC failed"
$wsOpen.Cells.Item(10, 2).Value = 300
$wsOpen.Cells.Item(11, 1).Value = "General technical question"
$wsOpen.Cells.Item(11, 2).Value = 120
$wsOpen.Cells.Item(12, 1).Value = "Error editing a variable"
$wsOpen.Cells.Item(12, 2).Value = 450
$wsOpen.Cells.Item(13, 1).Value = "No SAP data"
$wsOpen.Cells.Item(13, 2).Value = 200
$wsOpen.Cells.Item(14, 1).Value = "Scripting Problem"
$wsOpen.Cells.Item(14, 2).Value = 500
$wsOpen.Cells.Item(15, 1).Value = "Problems with editing fonts"
$wsOpen.Cells.Item(15, 2).Value = 250
$wsOpen.Cells.Item(16, 1).Value = "Authentication with JSON not possible"
$wsOpen.Cells.Item(16, 2).Value = 300
$wsOpen.Cells.Item(17, 1).Value = "Download error"
$wsOpen.Cells.Item(17, 2).Value = 300
$wsOpen.Cells.Item(18, 1).Value = "Special characters in password not accepted"
$wsOpen.Cells.Item(18, 2).Value = 100
$wsOpen.Cells.Item(19, 1).Value = "Error message"
$wsOpen.Cells.Item(19, 2).Value = 200

# Offene Tickets (German) row 20 - last new shared string
$wsOffene.Cells.Item(20, 1).Value = "Databinding Listview funktioniert nicht"
$wsOffene.Cells.Item(20, 2).Value = 250

# Offene Tickets (German) row 21 - reuses "SQL Async"
$wsOffene.Cells.Item(21, 1).Value = "SQL Async"
$wsOffene.Cells.Item(21, 2).Value = 100

# Give column A on "Offene Tickets" enough width to fit the new, longer labels
$wsOffene.Columns.Item(1).AutoFit()

# ---------------------------------------------------------------------------
# 2) Updated + extended statistics on "Solved tickets in a year"
# ---------------------------------------------------------------------------

$wsSolved.Cells.Item(2, 2).Value = 80
$wsSolved.Cells.Item(2, 3).Value = 92

$wsSolved.Cells.Item(3, 2).Value = 106
$wsSolved.Cells.Item(3, 3).Value = 102

$wsSolved.Cells.Item(4, 2).Value = 30
$wsSolved.Cells.Item(4, 3).Value = 30

$wsSolved.Cells.Item(5, 2).Value = 34
$wsSolved.Cells.Item(5, 3).Value = 40

$wsSolved.Cells.Item(6, 2).Value = 40
$wsSolved.Cells.Item(6, 3).Value = 54

$wsSolved.Cells.Item(7, 2).Value = 30
$wsSolved.Cells.Item(7, 3).Value = 40

$wsSolved.Cells.Item(8, 2).Value = 26
$wsSolved.Cells.Item(8, 3).Value = 30

$wsSolved.Cells.Item(9, 3).Value = 50

$wsSolved.Cells.Item(10, 2).Value = 60
$wsSolved.Cells.Item(10, 3).Value = 70

$wsSolved.Cells.Item(11, 2).Value = 10
$wsSolved.Cells.Item(11, 3).Value = 20

$wsSolved.Cells.Item(12, 2).Value = 20
$wsSolved.Cells.Item(12, 3).Value = 20

$wsSolved.Cells.Item(13, 2).Value = 16
$wsSolved.Cells.Item(13, 3).Value = 20

# New rows 14-20, column A carries the same date-number-format as the rows above it
function Add-SolvedRow($row, $dateVal, $b, $c) {
    $wsSolved.Cells.Item(13, 1).Copy() | Out-Null
    $wsSolved.Cells.Item($row, 1).PasteSpecial(-4122) | Out-Null
    $wsSolved.Cells.Item($row, 1).Value = $dateVal
    $wsSolved.Cells.Item($row, 2).Value = $b
    $wsSolved.Cells.Item($row, 3).Value = $c
}

Add-SolvedRow 14 44105 40 46
Add-SolvedRow 15 44136 50 56
Add-SolvedRow 16 44166 20 40
Add-SolvedRow 17 44197 28 30
Add-SolvedRow 18 44228 40 52
Add-SolvedRow 19 44256 20 25
Add-SolvedRow 20 44287 16 24

# ---------------------------------------------------------------------------
# 3) Remove the "Geschlossene Tickets pro Jahr" sheet (duplicate of
#    "Solved tickets in a year")
# ---------------------------------------------------------------------------

$wsGesch.Delete()

# ---------------------------------------------------------------------------
# 4) Selection / active-sheet bookkeeping to match the final view state
# ---------------------------------------------------------------------------

$wsOpen.Range("C20").Select()
$wsOffene.Range("A21").Select()

$wsSolved.Activate()
$wsSolved.Range("C13").Select()
